$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Index")
$ws.Range("A1").Value = "Test"
